$d = $word.ActiveDocument

# --- Change 1: mark the run holding the inline picture (the one with
#     <w:lastRenderedPageBreak/> + <w:drawing>) as NoProofing, which
#     serializes to <w:rPr><w:noProof/></w:rPr> on that run. ---
$shape = $d.InlineShapes.Item(1)
$shape.Range.NoProofing = $true

# --- Change 2: highlight (bright green) the whole paragraph
#     "Load the images from the folder if it was created beforehand (PRIO:1)"
#     including its paragraph mark, so both runs AND the paragraph mark's
#     rPr (w:pPr/w:rPr) get <w:highlight w:val="green"/>. ---
$find = $d.Content
$found = $find.Find.Execute(
    "Load the images from the folder if it was created beforehand (PRIO:1)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $find.Paragraphs(1).Range
    $para.Select()
    $word.Selection.Font.HighlightColorIndex = 4
}
